# "añadidos 125 y 250 hasta el año 2003"
# Adds a new "Class" column to the Table_2 (Equipos sheet / queryTable-backed
# ListObject) and fills every data row with "motogp".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipos")

# Expand the query table (Table_2) with a new, unbound column.
$lo = $ws.ListObjects.Item("Table_2")
$newCol = $lo.ListColumns.Add()

# Header (this also renames the new ListColumn to "Class").
$ws.Range("H1").Value = "Class"

# Fill the new column's data body with "motogp" for every existing row.
$ws.Range("H2:H26").Value = "motogp"

# Pick up the distinct "General" number-format style the real column carries.
$ws.Range("H2:H26").NumberFormat = "General"

# Restore the sheet selection to the newly added column, as left by the user.
$ws.Activate()
$ws.Range("H2:H26").Select()
